$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("attributes")

# Insert a new row at position 3 - this shifts the existing rows 3 (sex),
# 4 (ethnicity) and 5 (auto_id) down to rows 4, 5 and 6 respectively.
$ws.Rows.Item(3).Insert()

# New row 3: experimentID attribute definition
$ws.Cells.Item(3, 1).Value = "experimentID"
$ws.Cells.Item(3, 2).Value = "ExperimentID"
$ws.Cells.Item(3, 3).Value = "experiment ID"
$ws.Cells.Item(3, 4).Value = "rd3_portal_demographics"
$ws.Cells.Item(3, 5).Value = "string"
$ws.Cells.Item(3, 6).Value = $false
$ws.Cells.Item(3, 7).Value = $false
$ws.Cells.Item(3, 8).Value = $false
$ws.Cells.Item(3, 9).Value = $true

# Row 4 (was "ethnicity" before the insert shifted it down): now the
# "sex" attribute, with refreshed label/description text.
$ws.Cells.Item(4, 1).Value = "sex"
$ws.Cells.Item(4, 2).Value = "Observed Sex"
$ws.Cells.Item(4, 3).Value = "Observed Sex"

# Row 5 (was "auto_id" before the insert shifted it down): now the
# "ethnicity" attribute, with refreshed label/description text.
$ws.Cells.Item(5, 1).Value = "ethnicity"
$ws.Cells.Item(5, 2).Value = "Calculated Ancestry"
$ws.Cells.Item(5, 3).Value = "Ancestry that was derived"
$ws.Cells.Item(5, 6).Value = $false
$ws.Cells.Item(5, 7).Value = $false
$ws.Cells.Item(5, 8).Value = $true

# Column widths for the attributes sheet (COM ColumnWidth applies a
# +5/6-character pixel-rounding pass, same as real Excel, so the inputs
# here are pre-compensated to land on the intended final widths).
$ws.Columns.Item(1).ColumnWidth = 11.330729166666666
$ws.Columns.Item(2).ColumnWidth = 22.166666666666668
$ws.Columns.Item(3).ColumnWidth = 24.166666666666668
$ws.Columns.Item(4).ColumnWidth = 21.666666666666668

# Update the active selection to match
$ws.Range("A5").Select() | Out-Null
